$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values to replicate (same as existing row 111) for new rows 112-146
$values = @("KAGS #7158", 1.01, 57, "Jett", 23.3, 14, 138, 156, 0, 15, "Nickel", "['Rusher', 'Straight Up Winner']")

for ($r = 112; $r -le 146; $r++) {
    for ($c = 1; $c -le 12; $c++) {
        $ws.Cells.Item($r, $c).Value = $values[$c - 1]
    }
}
